$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: add new columns E, F, G (copy header style from B1) ---
$ws.Range("B1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$ws.Range("E1").Value = "C/A Lag"
$ws.Range("F1").Value = "LF Lag"
$ws.Range("G1").Value = "FFR Lag"

# --- Row 2: C/A Lag ---
$ws.Range("C2").Value = "'-0.065"
$ws.Range("D2").Value = "-0.036***"
$ws.Range("E2").Value = "1.0***"
$ws.Range("F2").Value = "'-0.0"
$ws.Range("G2").Value = "'-0.0"

# --- Row 3: LF Lag ---
$ws.Range("C3").Value = "-0.265***"
$ws.Range("D3").Value = "'-0.003"
$ws.Range("E3").Value = "-0.0**"
$ws.Range("F3").Value = "1.0***"
$ws.Range("G3").Value = "-0.0***"

# --- Row 4: FFR Lag ---
$ws.Range("C4").Value = "6.287***"
$ws.Range("D4").Value = "0.283***"
$ws.Range("E4").Value = "0.0***"
$ws.Range("F4").Value = "'0.0"
$ws.Range("G4").Value = "1.0***"

# --- Row 5: Constant ---
$ws.Range("C5").Value = "'-0.431"
$ws.Range("D5").Value = "'-0.13"
$ws.Range("E5").Value = "0.0*"
$ws.Range("F5").Value = "'0.0"
$ws.Range("G5").Value = "'0.0"

# --- Row 6: r2_adj (numeric values) ---
$ws.Range("C6").Value = 0.68
$ws.Range("D6").Value = 0.19
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1

# --- Strip the "quote prefix" formatting that got applied to the
#     apostrophe-escaped numeric-looking text cells above, so that no
#     extra style ends up attached to those (unstyled) data cells. ---
foreach ($addr in @("C2","F2","G2","D3","F4","C5","D5","F5","G5")) {
    $ws.Range($addr).Style = "Normal"
}
